$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (rows 5-9): only the "nb" counts in column B change ---
$ws.Range("B5").Value = 133530
$ws.Range("B6").Value = 30266
$ws.Range("B7").Value = 9171
$ws.Range("B8").Value = 3432
$ws.Range("B9").Value = 6454

# --- Table 2 (rows 13-17): both B (nb) and C (Kf) columns change ---
$ws.Range("B13").Value = 133530
$ws.Range("C13").Value = 0

$ws.Range("B14").Value = 30266
$ws.Range("C14").Value = 0

$ws.Range("B15").Value = 9170
$ws.Range("C15").Value = 6

$ws.Range("B16").Value = 3428
$ws.Range("C16").Value = 10

$ws.Range("B17").Value = 6451
$ws.Range("C17").Value = 6

# H18's total now needs to include row 16 as well (SUM(H12:H16) instead of SUM(H12:H15))
$ws.Range("H18").Formula = "=SUM(H12:H16)"

# --- Table 3 (rows 22-26): both B (nb) and C (Kf) columns change ---
$ws.Range("B22").Value = 133530
$ws.Range("C22").Value = 0

$ws.Range("B23").Value = 30266
$ws.Range("C23").Value = 1

$ws.Range("B24").Value = 9171
$ws.Range("C24").Value = 2467

$ws.Range("B25").Value = 966
$ws.Range("C25").Value = 2466

$ws.Range("B26").Value = 6454
$ws.Range("C26").Value = 0
